# Add NFL 2023 weeks 15-16 games (rows 482-513), trailer rows (514-527),
# matching the "Add files via upload" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("games")

# season, week, date(serial), team1, team2, score1, score2  -- keyed by row
$data = @(
  @(482, 15, 45274, "LAC", "LV",  21, 63),
  @(483, 15, 45276, "MIN", "CIN", 24, 27),
  @(484, 15, 45276, "PIT", "IND", 13, 30),
  @(485, 15, 45276, "DEN", "DET", 17, 42),
  @(486, 15, 45277, "ATL", "CAR", 7,  9),
  @(487, 15, 45277, "CHI", "CLE", 17, 20),
  @(488, 15, 45277, "NYG", "NO",  6,  24),
  @(489, 15, 45277, "TB",  "GB",  34, 20),
  @(490, 15, 45277, "HOU", "TEN", 19, 16),
  @(491, 15, 45277, "NYJ", "MIA", 0,  30),
  @(492, 15, 45277, "KC",  "NE",  27, 17),
  @(493, 15, 45277, "WAS", "LA",  20, 28),
  @(494, 15, 45277, "SF",  "ARI", 45, 29),
  @(495, 15, 45277, "DAL", "BUF", 10, 31),
  @(496, 15, 45277, "BAL", "JAX", 23, 7),
  @(497, 15, 45278, "PHI", "SEA", 17, 20),
  @(498, 16, 45281, "NO",  "LA ", 22, 30),
  @(499, 16, 45283, "CIN", "PIT", 11, 34),
  @(500, 16, 45283, "BUF", "LAC", 24, 22),
  @(501, 16, 45284, "CLE", "HOU", 36, 22),
  @(502, 16, 45284, "DET", "MIN", 30, 24),
  @(503, 16, 45284, "SEA", "TEN", 20, 17),
  @(504, 16, 45284, "IND", "ATL", 10, 29),
  @(505, 16, 45284, "WAS", "NYJ", 28, 30),
  @(506, 16, 45284, "GB",  "CAR", 33, 30),
  @(507, 16, 45284, "JAX", "TB",  12, 30),
  @(508, 16, 45284, "ARI", "CHI", 16, 27),
  @(509, 16, 45284, "DAL", "MIA", 20, 22),
  @(510, 16, 45284, "NE",  "DEN", 26, 23),
  @(511, 16, 45285, "LV",  "KC",  20, 14),
  @(512, 16, 45285, "NYG", "PHI", 25, 33),
  @(513, 16, 45285, "BAL", "SF",  33, 19)
)

foreach ($row in $data) {
  $r = $row[0]
  $ws.Cells.Item($r, 1).Value = 2023
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
  $ws.Cells.Item($r, 5).Value = $row[4]
  $ws.Cells.Item($r, 6).Value = $row[5]
  $ws.Cells.Item($r, 7).Value = $row[6]
}

# H482:H515 all carry the same relative formula (home_team = team2).
# Split the range write at the pre-existing shared-formula boundary
# (H483:H501 was already one shared group) so the new cells below it
# start a clean shared group of their own instead of corrupting the
# existing one.
$ws.Range("H482:H501").Formula = "=E482"
$ws.Range("H502:H515").Formula = "=E502"

# Give the new date cells (C482:C514) the same number format as the
# existing date column (copy format only, reuse the existing style).
$ws.Cells.Item(481, 3).Copy()
$ws.Range("C482:C514").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 514: season only (no game played), date cell stays styled/blank.
$ws.Cells.Item(514, 1).Value = 2023
$ws.Cells.Item(514, 3).ClearContents()

# Row 515: season only, no date styling.
$ws.Cells.Item(515, 1).Value = 2023

# Rows 516-527: season-only trailer rows.
for ($r = 516; $r -le 527; $r++) {
  $ws.Cells.Item($r, 1).Value = 2023
}

# Restore the view to what the commit shows (scrolled down / new selection).
$ws.Activate()
$ws.Range("B514").Select()
